# Update Sheets via scheduled runner: refresh market price / profit figures
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 8846.846
$ws.Cells.Item(69, 9).Value = 8749.5
$ws.Cells.Item(69, 11).Value = 26248.5
$ws.Cells.Item(69, 13).Value = -25374.5
$ws.Cells.Item(70, 8).Value = 8510.444
$ws.Cells.Item(70, 9).Value = 2630.6924
$ws.Cells.Item(70, 10).Value = 13970.214
$ws.Cells.Item(70, 11).Value = 7892.0772
$ws.Cells.Item(70, 12).Value = 41910.642
$ws.Cells.Item(70, 13).Value = -7622.0772
$ws.Cells.Item(70, 14).Value = -42450.642
$ws.Cells.Item(72, 8).Value = 8846.846
$ws.Cells.Item(72, 9).Value = 8749.5
$ws.Cells.Item(72, 11).Value = 78745.5
$ws.Cells.Item(72, 13).Value = -74377.5
$ws.Cells.Item(73, 8).Value = 8510.444
$ws.Cells.Item(73, 9).Value = 2630.6924
$ws.Cells.Item(73, 10).Value = 13970.214
$ws.Cells.Item(73, 11).Value = 7892.0772
$ws.Cells.Item(73, 12).Value = 41910.642
$ws.Cells.Item(73, 13).Value = -6956.0772
$ws.Cells.Item(73, 14).Value = -43782.642
$ws.Cells.Item(88, 8).Value = 3475.5
$ws.Cells.Item(88, 9).Value = 3476.25
$ws.Cells.Item(88, 11).Value = 3476.25
$ws.Cells.Item(88, 13).Value = -3070.25
$ws.Cells.Item(91, 8).Value = 3475.5
$ws.Cells.Item(91, 9).Value = 3476.25
$ws.Cells.Item(91, 11).Value = 3476.25
$ws.Cells.Item(91, 13).Value = -2072.25
$ws.Cells.Item(92, 8).Value = 6252.0586
$ws.Cells.Item(92, 9).Value = 5415.4614
$ws.Cells.Item(92, 11).Value = 5415.4614
$ws.Cells.Item(92, 13).Value = -4167.4614
$ws.Cells.Item(96, 8).Value = 1620.5
$ws.Cells.Item(96, 9).Value = 1387.2
$ws.Cells.Item(96, 10).Value = 1853.8
$ws.Cells.Item(96, 11).Value = 4161.6
$ws.Cells.Item(96, 12).Value = 5561.4
$ws.Cells.Item(96, 13).Value = -2788.6
$ws.Cells.Item(96, 14).Value = -8307.4
$ws.Cells.Item(98, 8).Value = 3995.68
$ws.Cells.Item(98, 9).Value = 1887.0454
$ws.Cells.Item(98, 11).Value = 1887.0454
$ws.Cells.Item(98, 13).Value = -389.0454
$ws.Cells.Item(116, 8).Value = 11950.944
$ws.Cells.Item(116, 10).Value = 6524.143
$ws.Cells.Item(116, 12).Value = 6524.143
$ws.Cells.Item(116, 14).Value = -13408.143
$ws.Cells.Item(122, 8).Value = 3995.68
$ws.Cells.Item(122, 9).Value = 1887.0454
$ws.Cells.Item(122, 11).Value = 5661.1362
$ws.Cells.Item(122, 13).Value = -3211.1362
$ws.Cells.Item(127, 8).Value = 22604.715
$ws.Cells.Item(127, 10).Value = 2995
$ws.Cells.Item(127, 12).Value = 8985
$ws.Cells.Item(127, 14).Value = -18905
$ws.Cells.Item(135, 8).Value = 1462.0869
$ws.Cells.Item(135, 9).Value = 1451.45
$ws.Cells.Item(135, 11).Value = 13063.05
$ws.Cells.Item(135, 13).Value = -10528.05
$ws.Cells.Item(138, 8).Value = 2812.28
$ws.Cells.Item(138, 9).Value = 1702.1666
$ws.Cells.Item(138, 10).Value = 3288.043
$ws.Cells.Item(138, 11).Value = 5106.4998
$ws.Cells.Item(138, 12).Value = 9864.129000000001
$ws.Cells.Item(138, 13).Value = 33.5002000000004
$ws.Cells.Item(138, 14).Value = -20144.129

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8386.15
$ws.Cells.Item(2, 9).Value = 9169.666999999999
$ws.Cells.Item(2, 10).Value = 1334.5
$ws.Cells.Item(2, 11).Value = 9169.666999999999
$ws.Cells.Item(2, 12).Value = 1334.5
$ws.Cells.Item(2, 13).Value = -9056.666999999999
$ws.Cells.Item(61, 8).Value = 6421.1
$ws.Cells.Item(61, 9).Value = 5395.0527
$ws.Cells.Item(61, 11).Value = 5395.0527
$ws.Cells.Item(61, 13).Value = -5183.0527
$ws.Cells.Item(116, 8).Value = 8386.15
$ws.Cells.Item(116, 9).Value = 9169.666999999999
$ws.Cells.Item(116, 10).Value = 1334.5
$ws.Cells.Item(116, 11).Value = 9169.666999999999
$ws.Cells.Item(116, 12).Value = 1334.5
$ws.Cells.Item(116, 13).Value = -6875.666999999999
$ws.Cells.Item(136, 8).Value = 6421.1
$ws.Cells.Item(136, 9).Value = 5395.0527
$ws.Cells.Item(136, 11).Value = 16185.1581
$ws.Cells.Item(136, 13).Value = -13635.1581

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8386.15
$ws.Cells.Item(3, 9).Value = 9169.666999999999
$ws.Cells.Item(3, 10).Value = 1334.5
$ws.Cells.Item(3, 11).Value = 9169.666999999999
$ws.Cells.Item(3, 12).Value = 1334.5
$ws.Cells.Item(3, 13).Value = -9055.666999999999
$ws.Cells.Item(107, 8).Value = 2911.353
$ws.Cells.Item(107, 9).Value = 2832.1667
$ws.Cells.Item(107, 11).Value = 2832.1667
$ws.Cells.Item(107, 13).Value = -912.1667000000002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 599865.5600000001
$ws.Cells.Item(31, 9).Value = 1430807.8
$ws.Cells.Item(31, 10).Value = 71084.17999999999
$ws.Cells.Item(31, 11).Value = 1430807.8
$ws.Cells.Item(31, 12).Value = 71084.17999999999
$ws.Cells.Item(31, 13).Value = -1430512.8
$ws.Cells.Item(31, 14).Value = -71674.17999999999
$ws.Cells.Item(34, 8).Value = 599865.5600000001
$ws.Cells.Item(34, 9).Value = 1430807.8
$ws.Cells.Item(34, 10).Value = 71084.17999999999
$ws.Cells.Item(34, 11).Value = 1430807.8
$ws.Cells.Item(34, 12).Value = 71084.17999999999
$ws.Cells.Item(34, 13).Value = -1430605.8
$ws.Cells.Item(34, 14).Value = -71488.17999999999
$ws.Cells.Item(62, 8).Value = 7165
$ws.Cells.Item(62, 9).Value = 7165
$ws.Cells.Item(62, 11).Value = 7165
$ws.Cells.Item(62, 13).Value = -6541
$ws.Cells.Item(65, 8).Value = 7165
$ws.Cells.Item(65, 9).Value = 7165
$ws.Cells.Item(65, 11).Value = 35825
$ws.Cells.Item(65, 13).Value = -32705
$ws.Cells.Item(107, 8).Value = 4785.5527
$ws.Cells.Item(107, 9).Value = 872.875
$ws.Cells.Item(107, 11).Value = 872.875
$ws.Cells.Item(107, 13).Value = 1047.125
$ws.Cells.Item(132, 8).Value = 6109.5264
$ws.Cells.Item(132, 9).Value = 3839.3572
$ws.Cells.Item(132, 11).Value = 11518.0716
$ws.Cells.Item(132, 13).Value = -8988.071599999999
$ws.Cells.Item(134, 8).Value = 6308.5
$ws.Cells.Item(134, 9).Value = 6708
$ws.Cells.Item(134, 10).Value = 3512
$ws.Cells.Item(134, 11).Value = 20124
$ws.Cells.Item(134, 12).Value = 10536
$ws.Cells.Item(134, 13).Value = -17589
$ws.Cells.Item(134, 14).Value = -15606

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 2199.5
$ws.Cells.Item(69, 10).Value = 2399
$ws.Cells.Item(69, 12).Value = 7197
$ws.Cells.Item(69, 14).Value = -8819
$ws.Cells.Item(72, 8).Value = 2199.5
$ws.Cells.Item(72, 10).Value = 2399
$ws.Cells.Item(72, 12).Value = 21591
$ws.Cells.Item(72, 14).Value = -29703
$ws.Cells.Item(107, 8).Value = 518.5599999999999
$ws.Cells.Item(107, 9).Value = 528.6667
$ws.Cells.Item(107, 10).Value = 509.23077
$ws.Cells.Item(107, 11).Value = 1586.0001
$ws.Cells.Item(107, 12).Value = 1527.69231
$ws.Cells.Item(107, 13).Value = 333.9999
$ws.Cells.Item(107, 14).Value = -5367.69231
$ws.Cells.Item(121, 8).Value = 2527.25
$ws.Cells.Item(121, 9).Value = 1408.625
$ws.Cells.Item(121, 10).Value = 3645.875
$ws.Cells.Item(121, 11).Value = 4225.875
$ws.Cells.Item(121, 12).Value = 10937.625
$ws.Cells.Item(121, 13).Value = -2915.875
$ws.Cells.Item(121, 14).Value = -13557.625
$ws.Cells.Item(132, 8).Value = 9786.846
$ws.Cells.Item(132, 10).Value = 2200
$ws.Cells.Item(132, 12).Value = 19800
$ws.Cells.Item(132, 14).Value = -24860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2712.3103
$ws.Cells.Item(80, 9).Value = 2664.32
$ws.Cells.Item(80, 11).Value = 2664.32
$ws.Cells.Item(80, 13).Value = -1666.32
$ws.Cells.Item(83, 8).Value = 2712.3103
$ws.Cells.Item(83, 9).Value = 2664.32
$ws.Cells.Item(83, 11).Value = 13321.6
$ws.Cells.Item(83, 13).Value = -8329.6
$ws.Cells.Item(104, 8).Value = 32450
$ws.Cells.Item(104, 10).Value = 32450
$ws.Cells.Item(104, 12).Value = 32450
$ws.Cells.Item(104, 14).Value = -39438
$ws.Cells.Item(126, 8).Value = 7478.16
$ws.Cells.Item(126, 9).Value = 6938.2
$ws.Cells.Item(126, 10).Value = 7838.1333
$ws.Cells.Item(126, 11).Value = 20814.6
$ws.Cells.Item(126, 12).Value = 23514.3999
$ws.Cells.Item(126, 13).Value = -18344.6
$ws.Cells.Item(126, 14).Value = -28454.3999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 21967.637
$ws.Cells.Item(40, 9).Value = 21967.637
$ws.Cells.Item(40, 11).Value = 21967.637
$ws.Cells.Item(40, 13).Value = -21831.637
$ws.Cells.Item(68, 8).Value = 2160.375
$ws.Cells.Item(68, 9).Value = 2144.4333
$ws.Cells.Item(68, 10).Value = 2399.5
$ws.Cells.Item(68, 11).Value = 2144.4333
$ws.Cells.Item(68, 12).Value = 2399.5
$ws.Cells.Item(68, 13).Value = -1395.4333
$ws.Cells.Item(68, 14).Value = -3897.5
$ws.Cells.Item(71, 8).Value = 2160.375
$ws.Cells.Item(71, 9).Value = 2144.4333
$ws.Cells.Item(71, 10).Value = 2399.5
$ws.Cells.Item(71, 11).Value = 10722.1665
$ws.Cells.Item(71, 12).Value = 11997.5
$ws.Cells.Item(71, 13).Value = -6978.166500000001
$ws.Cells.Item(71, 14).Value = -19485.5
$ws.Cells.Item(82, 8).Value = 1310.2
$ws.Cells.Item(82, 9).Value = 1197.762
$ws.Cells.Item(82, 10).Value = 1572.5555
$ws.Cells.Item(82, 11).Value = 1197.762
$ws.Cells.Item(82, 12).Value = 1572.5555
$ws.Cells.Item(82, 13).Value = -836.7619999999999
$ws.Cells.Item(82, 14).Value = -2294.5555
$ws.Cells.Item(85, 8).Value = 1310.2
$ws.Cells.Item(85, 9).Value = 1197.762
$ws.Cells.Item(85, 10).Value = 1572.5555
$ws.Cells.Item(85, 11).Value = 1197.762
$ws.Cells.Item(85, 12).Value = 1572.5555
$ws.Cells.Item(85, 13).Value = 50.23800000000006
$ws.Cells.Item(85, 14).Value = -4068.5555
$ws.Cells.Item(93, 8).Value = 1166.5555
$ws.Cells.Item(93, 9).Value = 1217.9333
$ws.Cells.Item(93, 10).Value = 909.6667
$ws.Cells.Item(93, 11).Value = 1217.9333
$ws.Cells.Item(93, 12).Value = 909.6667
$ws.Cells.Item(93, 13).Value = 30.06670000000008
$ws.Cells.Item(93, 14).Value = -3405.6667
$ws.Cells.Item(132, 8).Value = 9897.25
$ws.Cells.Item(132, 9).Value = 8468.299999999999
$ws.Cells.Item(132, 10).Value = 10691.111
$ws.Cells.Item(132, 11).Value = 25404.9
$ws.Cells.Item(132, 12).Value = 32073.333
$ws.Cells.Item(132, 13).Value = -22874.9
$ws.Cells.Item(132, 14).Value = -37133.333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10568.625
$ws.Cells.Item(62, 9).Value = 4789.5
$ws.Cells.Item(62, 11).Value = 4789.5
$ws.Cells.Item(62, 13).Value = -4165.5
$ws.Cells.Item(65, 8).Value = 10568.625
$ws.Cells.Item(65, 9).Value = 4789.5
$ws.Cells.Item(65, 11).Value = 23947.5
$ws.Cells.Item(65, 13).Value = -20827.5
$ws.Cells.Item(107, 8).Value = 1005.0769
$ws.Cells.Item(107, 9).Value = 956
$ws.Cells.Item(107, 11).Value = 2868
$ws.Cells.Item(107, 13).Value = -948
